# Latest in progress commit refactor of tests.xlsx:
#  - "suite-demo1" sheet had a duplicated block of rows (2-8 duplicated as
#    12-18). The duplicate block (rows 12-18) is removed.
#  - The "open" row's Target value changes from "/" to "/domainname/".
#  - Selection/cursor markers are updated to reflect the new, smaller sheet
#    (C4 on suite-demo1, B6 on "data set 1").

$wb = $excel.ActiveWorkbook

# --- "suite-demo1" sheet: drop the duplicated second block, fix the
#     domain-open target value, and move the active cell onto the
#     remaining data. ---
$ws2 = $wb.Worksheets.Item("suite-demo1")
$ws2.Activate()

# Remove the duplicate block of rows (rows 12 through 18).
[void]$ws2.Rows("12:18").Delete()

# "open" row's Target cell: "/" -> "/domainname/"
$ws2.Range("C3").Value = "/domainname/"

# Cursor ends up on C4 after the edit.
[void]$ws2.Range("C4").Select()

# --- "data set 1" sheet: only the selection/active cell moved (to B6). ---
$ws3 = $wb.Worksheets.Item("data set 1")
$ws3.Activate()
[void]$ws3.Range("B6").Select()

# Leave the originally active/tab-selected sheet ("suite-demo1") active.
$ws2.Activate()
